# regen sval data to filter save games
# Replaces the B:E and G (computed sum) values for rows 2-21 on Sheet1
# with the regenerated statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.01253208636536152;  C = 0.002658071450198252;  D = 0.1496068669990043;  E = 0.5333859586016987; G = 0.6981829834162627 }
    3  = @{ B = 3.272327238179451;    C = 1.626987699542094;     D = 0.7210945179870265;  E = 0.5333859586016987; G = 6.15379541431027 }
    4  = @{ B = 0.6545652718822623;   C = 1.626987699542094;     D = 3.223369029078222;   E = 0.5333859586016987; G = 6.038307959104277 }
    5  = @{ B = 0.6545652718822623;   C = 0.3048912486333797;    D = 3.223369029078222;   E = 0.5333859586016987; G = 4.716211508195562 }
    6  = @{ B = 0.003078177322033415; C = 0.002658071450198252;  D = 0.7210945179870265;  E = 0.5333859586016987; G = 1.260216725360957 }
    7  = @{ B = 3.272327238179451;    C = 1.626987699542094;     D = 0.7210945179870265;  E = 0.5333859586016987; G = 6.15379541431027 }
    8  = @{ B = 3.272327238179451;    C = 1.626987699542094;     D = 3.223369029078222;   E = 0.5333859586016987; G = 8.656069925401464 }
    9  = @{ B = 3.272327238179451;    C = 1.626987699542094;     D = 0.7210945179870265;  E = 0.5333859586016987; G = 6.15379541431027 }
    10 = @{ B = 0.6545652718822623;   C = 1.626987699542094;     D = 3.223369029078222;   E = 0.5333859586016987; G = 6.038307959104277 }
    11 = @{ B = 1.445647641019636;    C = 1.626987699542094;     D = 0.7210945179870265;  E = 0.5333859586016987; G = 4.327115817150455 }
    12 = @{ B = 0.1169995834814548;   C = 0.3048912486333797;    D = 0.1496068669990043;  E = 0.5333859586016987; G = 1.104883657715537 }
    13 = @{ B = 3.272327238179451;    C = 1.626987699542094;     D = 3.223369029078222;   E = 0.5333859586016987; G = 8.656069925401464 }
    14 = @{ B = 3.272327238179451;    C = 1.626987699542094;     D = 0.1496068669990043;  E = 0.5333859586016987; G = 5.582307763322248 }
    15 = @{ B = 3.272327238179451;    C = 1.626987699542094;     D = 0.7210945179870265;  E = 0.5333859586016987; G = 6.15379541431027 }
    16 = @{ B = 0.003078177322033415; C = 0.002658071450198252;  D = 0.1496068669990043;  E = 0.5333859586016987; G = 0.6887290743729346 }
    17 = @{ B = 3.272327238179451;    C = 1.626987699542094;     D = 0.1496068669990043;  E = 0.5333859586016987; G = 5.582307763322248 }
    18 = @{ B = 0.2881169905109251;   C = 0.3048912486333797;    D = 0.7210945179870265;  E = 0.5333859586016987; G = 1.84748871573303 }
    19 = @{ B = 3.272327238179451;    C = 1.626987699542094;     D = 0.1496068669990043;  E = 0.5333859586016987; G = 5.582307763322248 }
    20 = @{ B = 1.445647641019636;    C = 1.626987699542094;     D = 0.7210945179870265;  E = 0.5333859586016987; G = 4.327115817150455 }
    21 = @{ B = 1.445647641019636;    C = 1.626987699542094;     D = 0.1496068669990043;  E = 0.5333859586016987; G = 3.755628166162433 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
